$wb = $excel.ActiveWorkbook

# --- Existing "Tables" sheet: add the two table rows ---
$tables = $wb.Worksheets.Item("Tables")

$tables.Cells.Item(2,1).Value = "Table 1"
$c = $tables.Cells.Item(2,2)
$c.NumberFormat = "@"
$c.Value = "1"
$tables.Cells.Item(2,3).Value = $false

$tables.Cells.Item(3,1).Value = "Table 2"
$c = $tables.Cells.Item(3,2)
$c.NumberFormat = "@"
$c.Value = "2"
$tables.Cells.Item(3,3).Value = $false

# --- New sheets: "Table 1" and "Table 2", each with a menu header row ---
$t1 = $wb.Worksheets.Add($null, $tables)
$t1.Name = "Table 1"
$t1.Cells.Item(1,1).Value = "SNo"
$t1.Cells.Item(1,2).Value = "Dish"
$t1.Cells.Item(1,3).Value = "Quantity"

$t2 = $wb.Worksheets.Add($null, $t1)
$t2.Name = "Table 2"
$t2.Cells.Item(1,1).Value = "SNo"
$t2.Cells.Item(1,2).Value = "Dish"
$t2.Cells.Item(1,3).Value = "Quantity"

$tables.Activate()
